# "ema trader 3 - dev"
# Rename the short/long/additional_ema_len parameter labels to the new
# sema/lema/slema naming, and update their corresponding values (column B),
# letting the dependent formulas in column C recalculate automatically.
# Finally move the active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "sema"
$ws.Range("A5").Value = "lema"
$ws.Range("A6").Value = "slema"

$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 100
$ws.Range("B6").Value = 25

$ws.Range("A1").Select()
